$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: set the "Last Invoice Date" (column D) which was previously blank ---
# Copy formatting from an existing date cell (D2) so the numeric/date style (s="12")
# is reused instead of creating a brand-new style entry, then set the date value.
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 45971

# --- Insert a new row at position 27 (old rows 27/28 shift down to 28/29) ---
$ws.Rows("27").Insert()

# Re-apply the same row height used by all the other data rows.
$ws.Rows("27").RowHeight = 13.05

# Reuse the "blank date" style (s="13", from D6) for the new row's D cell instead
# of the date style that Insert() copied down from the row above.
$ws.Range("D6").Copy()
$ws.Range("D27").PasteSpecial(-4122)

# Reuse the formatting of the (empty, styleless) F column cell so the new row gets
# a matching placeholder cell in column F.
$ws.Range("F26").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# --- Fill in the new row's data ---
$ws.Range("A27").Value = "NATURE PATHWAYS ELC"
$ws.Range("B27").Value = "Zigan, Gerald L"
$ws.Range("C27").Value = "015"
$ws.Range("E27").Value = "0008366"
